$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Chinese translation for EFFECT_TEMP_ITEM (row 9, column C)
$ws.Range("C9").Value = "在接下来的{0}波敌袭中"

# Add new row 10: ATA_ITEM_CANNED_FOOD translation
$ws.Range("A10").Value = "ATA_ITEM_CANNED_FOOD"
$ws.Range("B10").Value = "Canned Food"
$ws.Range("C10").Value = "罐装食品"

# Add new row 11: EFFECT_GAIN_EVERY_CONSUMABLE translation
$ws.Range("A11").Value = "EFFECT_GAIN_EVERY_CONSUMABLE"
$ws.Range("B11").Value = "Gain 1 {0} when picking up every {1} consumables"
$ws.Range("C11").Value = "每当你拾起{0}个消耗品时获得1{1}"

# Update column widths to match new content
# (values chosen so the runtime's pixel-quantized stored width lands as close
# as possible to the authored widths of 31.453125 / 43.90625 / 40.54296875)
$ws.Columns.Item(1).ColumnWidth = 30.67
$ws.Columns.Item(2).ColumnWidth = 43
$ws.Columns.Item(3).ColumnWidth = 39.67

# Update the active selection to B10, matching the author's final cursor position
$ws.Range("B10").Select() | Out-Null
